# Commit: "Fruta / hortaliza, semanal" — weekly update adds one new
# market-price observation. It lands as a new row 103 in the data table
# (pushing the existing rows 103-149 down to 104-150); the rest of the
# sheet is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 103, shifting rows 103:149 down to 104:150.
$ws.Rows("103:103").Insert()

# Populate the new row 103 with the new weekly observation. Columns that
# are constant across this whole sub-sheet (Mercado ID, Mercado, Región,
# Codreg, Categoría ID/Categoría/Variedad/Calidad, Unidad de
# comercialización, Kg o Unidades, Clasificación) keep the same values as
# every other row.
$row = 103

$ws.Cells.Item($row, 1).Value2 = 7
$ws.Cells.Item($row, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value2 = "Ñuble"
$ws.Cells.Item($row, 4).Value2 = (Get-Date -Year 2021 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($row, 5).Value2 = 16
$ws.Cells.Item($row, 6).Value2 = 100112017
$ws.Cells.Item($row, 7).Value2 = "Apio"
$ws.Cells.Item($row, 8).Value2 = "Americana (o)"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 60
$ws.Cells.Item($row, 11).Value2 = 8000
$ws.Cells.Item($row, 12).Value2 = 8500
$ws.Cells.Item($row, 13).Value2 = 8250
$ws.Cells.Item($row, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item($row, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value2 = 1375
$ws.Cells.Item($row, 17).Value2 = 6
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date-formatted style the rest of
# column D uses (Insert() already copies the format from the row above,
# but set it explicitly to be safe).
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 1, 4).NumberFormat
